# Auto-generated edit script applying scheduled market-data refresh to leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 300.76923
$ws.Range("I4").Value = 273.63635
$ws.Range("J4").Value = 450
$ws.Range("K4").Value = 273.63635
$ws.Range("L4").Value = 450
$ws.Range("M4").Value = -159.63635
$ws.Range("N4").Value = -678

$ws.Range("H26").Value = 4633.3335
$ws.Range("I26").Value = 2666.6667
$ws.Range("J26").Value = 6600
$ws.Range("K26").Value = 2666.6667
$ws.Range("L26").Value = 6600
$ws.Range("M26").Value = -2322.6667
$ws.Range("N26").Value = -7288

$ws.Range("H32").Value = 4621.5557
$ws.Range("I32").Value = 4738.4
$ws.Range("J32").Value = 4475.5
$ws.Range("K32").Value = 4738.4
$ws.Range("L32").Value = 4475.5
$ws.Range("M32").Value = -4412.4
$ws.Range("N32").Value = -5127.5

$ws.Range("H42").Value = 12.666667
$ws.Range("I42").Value = 9
$ws.Range("J42").Value = 20
$ws.Range("K42").Value = 27
$ws.Range("L42").Value = 60
$ws.Range("M42").Value = 203
$ws.Range("N42").Value = -520

$ws.Range("H53").Value = 375.86667
$ws.Range("I53").Value = 217.85715
$ws.Range("J53").Value = 514.125
$ws.Range("K53").Value = 217.85715
$ws.Range("L53").Value = 514.125
$ws.Range("M53").Value = 419.14285

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H88").Value = 900
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 900
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 900
$ws.Range("N88").Value = -1712

$ws.Range("H91").Value = 900
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 900
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 900
$ws.Range("N91").Value = -3708

$ws.Range("H94").Value = 1206.2858
$ws.Range("I94").Value = 1206.2858
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1206.2858
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -755.2858000000001

$ws.Range("H132").Value = 938.5
$ws.Range("I132").Value = 938.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2815.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -285.5

$ws.Range("H138").Value = 6578.7617
$ws.Range("I138").Value = 4961.375
$ws.Range("J138").Value = 7574.077
$ws.Range("K138").Value = 14884.125
$ws.Range("L138").Value = 22722.231
$ws.Range("M138").Value = -9744.125
$ws.Range("N138").Value = -33002.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1099
$ws.Range("I61").Value = 1099
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1099
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -887

$ws.Range("H74").Value = 206
$ws.Range("I74").Value = 206
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 206
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 668

$ws.Range("H77").Value = 206
$ws.Range("I77").Value = 206
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 1030
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 3338

$ws.Range("H109").Value = 10000
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 10000
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 10000
$ws.Range("N109").Value = -12774

$ws.Range("H136").Value = 1099
$ws.Range("I136").Value = 1099
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3297
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -747

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 50001
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = 50001
$ws.Range("K108").Value = 0
$ws.Range("L108").Value = 50001
$ws.Range("N108").Value = -57681

$ws.Range("H122").Value = 1603
$ws.Range("I122").Value = 1395.6666
$ws.Range("J122").Value = 1914
$ws.Range("K122").Value = 4186.9998
$ws.Range("L122").Value = 5742
$ws.Range("M122").Value = -1736.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6000
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -6576

$ws.Range("H68").Value = 4166.5713
$ws.Range("I68").Value = 4443
$ws.Range("J68").Value = 3798
$ws.Range("K68").Value = 13329
$ws.Range("L68").Value = 11394
$ws.Range("M68").Value = -12518
$ws.Range("N68").Value = -13016

$ws.Range("H71").Value = 4166.5713
$ws.Range("I71").Value = 4443
$ws.Range("J71").Value = 3798
$ws.Range("K71").Value = 39987
$ws.Range("L71").Value = 34182
$ws.Range("M71").Value = -35931
$ws.Range("N71").Value = -42294

$ws.Range("H81").Value = 1500
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1500
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 4500
$ws.Range("N81").Value = -6746

$ws.Range("H84").Value = 1500
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1500
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 13500
$ws.Range("N84").Value = -24732

$ws.Range("H101").Value = 15000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 15000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 45000
$ws.Range("N101").Value = -49868

$ws.Range("H137").Value = 2712.2307
$ws.Range("I137").Value = 2404.75
$ws.Range("J137").Value = 2848.889
$ws.Range("K137").Value = 7214.25
$ws.Range("L137").Value = 8546.667000000001
$ws.Range("M137").Value = -2114.25
$ws.Range("N137").Value = -18746.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 5015000
$ws.Range("I35").Value = 5015000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 5015000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -5014702
$ws.Range("N35").ClearContents()

$ws.Range("H80").Value = 8715.076999999999
$ws.Range("I80").Value = 7613.7144
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 7613.7144
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -6615.7144

$ws.Range("H83").Value = 8715.076999999999
$ws.Range("I83").Value = 7613.7144
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 38068.572
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -33076.572

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 4299.2666
$ws.Range("I132").Value = 3449.3333
$ws.Range("J132").Value = 7699
$ws.Range("K132").Value = 10347.9999
$ws.Range("L132").Value = 23097
$ws.Range("M132").Value = -7817.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1958.7
$ws.Range("I7").Value = 2718.6
$ws.Range("J7").Value = 1198.8
$ws.Range("K7").Value = 2718.6
$ws.Range("L7").Value = 1198.8
$ws.Range("M7").Value = -2606.6
$ws.Range("N7").Value = -1422.8

$ws.Range("H38").Value = 21500
$ws.Range("I38").Value = 3000
$ws.Range("J38").Value = 40000
$ws.Range("K38").Value = 3000
$ws.Range("L38").Value = 40000
$ws.Range("M38").Value = -2590
$ws.Range("N38").Value = -40820

$ws.Range("H40").Value = 1599.8
$ws.Range("I40").Value = 1749.75
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 1749.75
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = -1613.75

$ws.Range("H55").Value = 1250
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1250
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 1250
$ws.Range("N55").Value = -1596

$ws.Range("H82").Value = 3750.6667
$ws.Range("I82").Value = 3002
$ws.Range("J82").Value = 4125
$ws.Range("K82").Value = 3002
$ws.Range("L82").Value = 4125
$ws.Range("M82").Value = -2641
$ws.Range("N82").Value = -4847

$ws.Range("H85").Value = 3750.6667
$ws.Range("I85").Value = 3002
$ws.Range("J85").Value = 4125
$ws.Range("K85").Value = 3002
$ws.Range("L85").Value = 4125
$ws.Range("M85").Value = -1754
$ws.Range("N85").Value = -6621

$ws.Range("H122").Value = 1544
$ws.Range("I122").Value = 1544
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4632
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2182

$ws.Range("H126").Value = 1958.7
$ws.Range("I126").Value = 2718.6
$ws.Range("J126").Value = 1198.8
$ws.Range("K126").Value = 8155.799999999999
$ws.Range("L126").Value = 3596.4
$ws.Range("M126").Value = -5685.799999999999
$ws.Range("N126").Value = -8536.4

$ws.Range("H132").Value = 4812
$ws.Range("I132").Value = 4765
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 14295
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -11765
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 1670100
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 2505000
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 2505000
$ws.Range("M18").Value = -127
$ws.Range("N18").Value = -2505346

$ws.Range("H29").Value = 15000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 15000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 15000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -15580

$ws.Range("H46").Value = 66666.336
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 66666.336
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 66666.336
$ws.Range("N46").Value = -67128.336

$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 4966.3335
$ws.Range("I122").Value = 4966.3335
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 14899.0005
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -12449.0005

$ws.Range("H126").Value = 835.7
$ws.Range("I126").Value = 918.875
$ws.Range("J126").Value = 503
$ws.Range("K126").Value = 2756.625
$ws.Range("L126").Value = 1509
$ws.Range("M126").Value = -286.625
$ws.Range("N126").Value = -6449

$ws.Range("H134").Value = 66666.336
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 66666.336
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 199999.008
$ws.Range("N134").Value = -205069.008
